# "fix akurasi yg salah" (fix wrong accuracy) — corrects the Reward values
# used to drive the three "My Bot" line charts on Sheet1 (columns D, H, L).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# My Bot 1 reward values (column D, rows 2-11)
$botOneRewards = 682.19, 382.61, 698.65, 747.69, 1491.12, 1910.56, 574.52, 1145.46, 1987.18, 907.07
for ($i = 0; $i -lt $botOneRewards.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $botOneRewards[$i]
}

# My Bot 2 reward values (column H, rows 2-11)
$botTwoRewards = 849, 1175, 1052.36, 1337.06, 1368.43, 414.36, 1155.45, 1286, 1056.45, 1262.81
for ($i = 0; $i -lt $botTwoRewards.Length; $i++) {
    $ws.Cells.Item($i + 2, 8).Value = $botTwoRewards[$i]
}

# My Bot 3 reward values (column L, rows 2-11)
$botThreeRewards = 1934.74, 2055.98, 1778.74, 1470.48, 1047.07, 1588.11, 2057.01, 1499.76, 937.93, 1588.05
for ($i = 0; $i -lt $botThreeRewards.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $botThreeRewards[$i]
}
